# Post-competition update: add features found from the text corpus.
# Append three new score rows (12-14) below the existing table on Sheet1,
# the first one labelled "addFeature" in column A, the other two being
# follow-up runs (no label, same metric columns B/C/D).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 - labelled result row
$ws.Cells.Item(12, 1).Value = "addFeature"
$ws.Cells.Item(12, 2).Value = 0.92893899999999996
$ws.Cells.Item(12, 3).Value = 0.92042000000000002
$ws.Cells.Item(12, 4).Value = 0.88973999999999998

# Row 13 - follow-up run, no label
$ws.Cells.Item(13, 2).Value = 0.92547449999999998
$ws.Cells.Item(13, 3).Value = 0.92284999999999995
$ws.Cells.Item(13, 4).Value = 0.89653000000000005

# Row 14 - follow-up run, no label
$ws.Cells.Item(14, 2).Value = 0.93502370000000001
$ws.Cells.Item(14, 3).Value = 0.92161999999999999
$ws.Cells.Item(14, 4).Value = 0.89295000000000002

# Leave the selection where the author left off editing
$ws.Range("B13").Select()
